$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.3312325492543664
$ws.Cells.Item(1, 2).Value = 0.33054250840825716
$ws.Cells.Item(2, 1).Value = -0.24093293693744755
$ws.Cells.Item(2, 2).Value = 0.23878311116014395
$ws.Cells.Item(3, 1).Value = -0.13583602361698865
$ws.Cells.Item(3, 2).Value = 0.13533224961756929
$ws.Cells.Item(4, 1).Value = -0.12333224978495849
$ws.Cells.Item(4, 2).Value = 0.12290416546665028
$ws.Cells.Item(5, 1).Value = -0.11690416609470677
$ws.Cells.Item(5, 2).Value = 0.11607005800634251
$ws.Cells.Item(6, 1).Value = -0.06952140940059115
$ws.Cells.Item(6, 2).Value = 0.06945875225206688
$ws.Cells.Item(7, 1).Value = -0.049458753007932899
$ws.Cells.Item(7, 2).Value = 0.049343590176309959
$ws.Cells.Item(8, 1).Value = -0.029343590937174646
$ws.Cells.Item(8, 2).Value = 0.029283937405797822
$ws.Cells.Item(9, 1).Value = -0.023283938059694087
$ws.Cells.Item(9, 2).Value = 0.0232439552021102
$ws.Cells.Item(10, 1).Value = -0.01724395585936378
$ws.Cells.Item(10, 2).Value = 0.01724448236458187
$ws.Cells.Item(11, 1).Value = -0.012744483009992535
$ws.Cells.Item(11, 2).Value = 0.012739487891707313
$ws.Cells.Item(12, 1).Value = -0.0067394885493308188
$ws.Cells.Item(12, 2).Value = 0.0067137646884498636
$ws.Cells.Item(13, 1).Value = -0.00071376534731282248
$ws.Cells.Item(13, 2).Value = 0.00070214737124452142
$ws.Cells.Item(14, 1).Value = 0.01129785192220556
$ws.Cells.Item(14, 2).Value = -0.011329186452171847
$ws.Cells.Item(15, 1).Value = 0.017329185794243251
$ws.Cells.Item(15, 2).Value = -0.017380290106411067
$ws.Cells.Item(16, 1).Value = 0.023380289450317449
$ws.Cells.Item(16, 2).Value = -0.02350068370277425
$ws.Cells.Item(17, 1).Value = 0.029500683051265852
$ws.Cells.Item(17, 2).Value = -0.029589319393379832
$ws.Cells.Item(18, 1).Value = -0.079609572101055903
$ws.Cells.Item(18, 2).Value = 0.079548210418323606
$ws.Cells.Item(19, 1).Value = -0.070548211041872477
$ws.Cells.Item(19, 2).Value = 0.070074575838816155
$ws.Cells.Item(20, 1).Value = -0.018013411165075155
$ws.Cells.Item(20, 2).Value = 0.01800425915150683
$ws.Cells.Item(21, 1).Value = -0.0090042597913786437
$ws.Cells.Item(21, 2).Value = 0.0089999993595428762
$ws.Cells.Item(22, 1).Value = -0.093947088182337879
$ws.Cells.Item(22, 2).Value = 0.0936344979318946
$ws.Cells.Item(23, 1).Value = -0.084634498576907191
$ws.Cells.Item(23, 2).Value = 0.084126712287129202
$ws.Cells.Item(24, 1).Value = -0.042126713204734045
$ws.Cells.Item(24, 2).Value = 0.041999999077262373
$ws.Cells.Item(25, 1).Value = -0.040548649248439261
$ws.Cells.Item(25, 2).Value = 0.040519564943608088
$ws.Cells.Item(26, 1).Value = -0.034519565586407452
$ws.Cells.Item(26, 2).Value = 0.034489983935721114
$ws.Cells.Item(27, 1).Value = -0.028489984579410432
$ws.Cells.Item(27, 2).Value = 0.028400873845886743
$ws.Cells.Item(28, 1).Value = -0.022400874493880174
$ws.Cells.Item(28, 2).Value = 0.022354246240972486
$ws.Cells.Item(29, 1).Value = -0.010354246939167311
$ws.Cells.Item(29, 2).Value = 0.010347252874440827
$ws.Cells.Item(30, 1).Value = 0.0096527463632307686
$ws.Cells.Item(30, 2).Value = -0.0098059969126098601
$ws.Cells.Item(31, 1).Value = 0.024805996191535229
$ws.Cells.Item(31, 2).Value = -0.024878785444911955
$ws.Cells.Item(32, 1).Value = -0.020616282984017964
$ws.Cells.Item(32, 2).Value = 0.020596381311274392

$ws.Columns.Item(1).ColumnWidth = 15.67
$ws.Columns.Item(2).ColumnWidth = 14.8
